# Update leve profit market-data columns (H:N) across all sheets.
# Values refreshed from the latest Universalis market snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H7").Value = 3770.8
$ws.Range("I7").Value = 3152.5
$ws.Range("J7").Value = 4183
$ws.Range("K7").Value = 3152.5
$ws.Range("L7").Value = 4183
$ws.Range("M7").Value = -3040.5
$ws.Range("N7").Value = -4407
$ws.Range("H9").Value = 73.5
$ws.Range("I9").Value = 73.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 73.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 95.5
$ws.Range("N9").Value = $null
$ws.Range("H14").Value = 3770.8
$ws.Range("I14").Value = 3152.5
$ws.Range("J14").Value = 4183
$ws.Range("K14").Value = 3152.5
$ws.Range("L14").Value = 4183
$ws.Range("M14").Value = -2961.5
$ws.Range("N14").Value = -4565
$ws.Range("H18").Value = 1278.25
$ws.Range("I18").Value = 1371.3334
$ws.Range("J18").Value = 999
$ws.Range("K18").Value = 1371.3334
$ws.Range("L18").Value = 999
$ws.Range("M18").Value = -1087.3334
$ws.Range("N18").Value = -1567
$ws.Range("H40").Value = 6413.5654
$ws.Range("J40").Value = 7599.5
$ws.Range("L40").Value = 7599.5
$ws.Range("N40").Value = -7949.5
$ws.Range("H64").Value = 9500
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = $null
$ws.Range("H67").Value = 9500
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = $null
$ws.Range("H92").Value = 288.75
$ws.Range("I92").Value = 288.75
$ws.Range("K92").Value = 288.75
$ws.Range("M92").Value = 959.25
$ws.Range("H106").Value = 7599.1
$ws.Range("I106").Value = 7331.6665
$ws.Range("K106").Value = 7331.6665
$ws.Range("M106").Value = -6700.6665
$ws.Range("H140").Value = 79260
$ws.Range("J140").Value = 79260
$ws.Range("L140").Value = 79260
$ws.Range("N140").Value = -89620

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H6").Value = 18001900
$ws.Range("I6").Value = 17145572
$ws.Range("K6").Value = 17145572
$ws.Range("M6").Value = -17145399
$ws.Range("H43").Value = 12499990
$ws.Range("J43").Value = 9999980
$ws.Range("L43").Value = 9999980
$ws.Range("N43").Value = -10000606
$ws.Range("H110").Value = 850
$ws.Range("I110").Value = 850
$ws.Range("K110").Value = 850
$ws.Range("M110").Value = 1195
$ws.Range("H122").Value = 1744.9231
$ws.Range("I122").Value = 1784
$ws.Range("J122").Value = 1614.6666
$ws.Range("K122").Value = 5352
$ws.Range("L122").Value = 4843.9998
$ws.Range("M122").Value = -2902
$ws.Range("N122").Value = -9743.9998

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H86").Value = 3835.8333
$ws.Range("I86").Value = 2390.3076
$ws.Range("K86").Value = 2390.3076
$ws.Range("M86").Value = -1267.3076
$ws.Range("H89").Value = 3835.8333
$ws.Range("I89").Value = 2390.3076
$ws.Range("K89").Value = 11951.538
$ws.Range("M89").Value = -6335.538
$ws.Range("H99").Value = 2621.875
$ws.Range("I99").Value = 1815.6
$ws.Range("K99").Value = 1815.6
$ws.Range("M99").Value = -317.5999999999999
$ws.Range("H105").Value = 1996
$ws.Range("I105").Value = 1997.5
$ws.Range("K105").Value = 1997.5
$ws.Range("M105").Value = -250.5
$ws.Range("H107").Value = 4422.696
$ws.Range("I107").Value = 3316.8
$ws.Range("K107").Value = 3316.8
$ws.Range("M107").Value = -1396.8

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H7").Value = 17141.834
$ws.Range("I7").Value = 25265.25
$ws.Range("J7").Value = 895
$ws.Range("K7").Value = 25265.25
$ws.Range("L7").Value = 895
$ws.Range("M7").Value = -25152.25
$ws.Range("N7").Value = -1121
$ws.Range("H22").Value = 2304.4
$ws.Range("I22").Value = 1769.5
$ws.Range("K22").Value = 1769.5
$ws.Range("M22").Value = -1419.5
$ws.Range("H33").Value = 1045.4286
$ws.Range("I33").Value = 1045.4286
$ws.Range("K33").Value = 1045.4286
$ws.Range("M33").Value = -666.4286
$ws.Range("H86").Value = 4900
$ws.Range("I86").Value = 4900
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4900
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3777
$ws.Range("N86").Value = $null
$ws.Range("H89").Value = 4900
$ws.Range("I89").Value = 4900
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 24500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -18884
$ws.Range("N89").Value = $null

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 1272
$ws.Range("I5").Value = 957.4
$ws.Range("J5").Value = 1796.3334
$ws.Range("K5").Value = 2872.2
$ws.Range("L5").Value = 5389.0002
$ws.Range("M5").Value = -2760.2
$ws.Range("N5").Value = -5613.0002
$ws.Range("H40").Value = 156.84616
$ws.Range("I40").Value = 23.625
$ws.Range("J40").Value = 370
$ws.Range("K40").Value = 94.5
$ws.Range("L40").Value = 1480
$ws.Range("M40").Value = -25.5
$ws.Range("N40").Value = -1618
$ws.Range("H75").Value = 4999
$ws.Range("J75").Value = 4999
$ws.Range("L75").Value = 14997
$ws.Range("N75").Value = -16993
$ws.Range("H78").Value = 4999
$ws.Range("J78").Value = 4999
$ws.Range("L78").Value = 44991
$ws.Range("N78").Value = -54975
$ws.Range("H135").Value = 1272
$ws.Range("I135").Value = 957.4
$ws.Range("J135").Value = 1796.3334
$ws.Range("K135").Value = 8616.6
$ws.Range("L135").Value = 16167.0006
$ws.Range("M135").Value = -6081.6
$ws.Range("N135").Value = -21237.0006

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H33").Value = 9998.375
$ws.Range("J33").Value = 9998.375
$ws.Range("L33").Value = 9998.375
$ws.Range("N33").Value = -10502.375
$ws.Range("H70").Value = 1158.1428
$ws.Range("I70").Value = 1017.8333
$ws.Range("K70").Value = 1017.8333
$ws.Range("M70").Value = -747.8333
$ws.Range("H73").Value = 1158.1428
$ws.Range("I73").Value = 1017.8333
$ws.Range("K73").Value = 1017.8333
$ws.Range("M73").Value = -81.83330000000001
$ws.Range("H102").Value = 2740.6667
$ws.Range("I102").Value = 2042
$ws.Range("J102").Value = 4836.6665
$ws.Range("K102").Value = 2042
$ws.Range("L102").Value = 4836.6665
$ws.Range("M102").Value = -420
$ws.Range("N102").Value = -8080.6665
$ws.Range("H122").Value = 3296.2222
$ws.Range("I122").Value = 2667.2856
$ws.Range("J122").Value = 5497.5
$ws.Range("K122").Value = 8001.8568
$ws.Range("L122").Value = 16492.5
$ws.Range("M122").Value = -5551.8568
$ws.Range("N122").Value = -21392.5
$ws.Range("H139").Value = 29999
$ws.Range("J139").Value = 29999
$ws.Range("L139").Value = 29999
$ws.Range("N139").Value = -40279

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H40").Value = 6100.722
$ws.Range("I40").Value = 6351.2856
$ws.Range("J40").Value = 5223.75
$ws.Range("K40").Value = 6351.2856
$ws.Range("L40").Value = 5223.75
$ws.Range("M40").Value = -6215.2856
$ws.Range("N40").Value = -5495.75
$ws.Range("H68").Value = 3992.3845
$ws.Range("I68").Value = 1920.2
$ws.Range("K68").Value = 1920.2
$ws.Range("M68").Value = -1171.2
$ws.Range("H71").Value = 3992.3845
$ws.Range("I71").Value = 1920.2
$ws.Range("K71").Value = 9601
$ws.Range("M71").Value = -5857
$ws.Range("H93").Value = 1310.5
$ws.Range("I93").Value = 1584
$ws.Range("K93").Value = 1584
$ws.Range("M93").Value = -336
$ws.Range("H136").Value = 3798
$ws.Range("I136").Value = 1395
$ws.Range("K136").Value = 4185
$ws.Range("M136").Value = -1635

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H41").Value = 35571.145
$ws.Range("J41").Value = 35553.6
$ws.Range("L41").Value = 35553.6
$ws.Range("N41").Value = -36333.6
$ws.Range("H45").Value = 21284.5
$ws.Range("I45").Value = 21284.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 21284.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -20793.5
$ws.Range("N45").Value = $null
$ws.Range("H74").Value = 9749
$ws.Range("J74").Value = 9749
$ws.Range("L74").Value = 9749
$ws.Range("N74").Value = -11621
$ws.Range("H77").Value = 9749
$ws.Range("J77").Value = 9749
$ws.Range("L77").Value = 29247
$ws.Range("N77").Value = -38607
$ws.Range("H96").Value = 1539.4
$ws.Range("I96").Value = 1566
$ws.Range("J96").Value = 1499.5
$ws.Range("K96").Value = 1566
$ws.Range("L96").Value = 1499.5
$ws.Range("M96").Value = -193
$ws.Range("N96").Value = -4245.5
$ws.Range("H100").Value = 1336
$ws.Range("I100").Value = 1037.6
$ws.Range("K100").Value = 2075.2
$ws.Range("M100").Value = -1534.2
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null
$ws.Range("H136").Value = 2401.8572
$ws.Range("I136").Value = 1739.4546
$ws.Range("J136").Value = 4830.6665
$ws.Range("K136").Value = 5218.3638
$ws.Range("L136").Value = 14491.9995
$ws.Range("M136").Value = -2668.3638
$ws.Range("N136").Value = -19591.9995
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
